$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) from 2023-11-13 (45243)
# to 2023-11-14 (45244) for the data rows (rows 2 through 10).
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 45244
}
